$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 260, shifting existing rows 260:378 down to 261:379
$ws.Rows("260").Insert()

# Populate the new row 260 with data (same catalog fields as the row that used to
# be at 260, but new observation values for D, L, M, N, O, P, S)
$ws.Range("A260").Value = 7
$ws.Range("B260").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C260").Value = "Ñuble"
$ws.Range("D260").Value = [DateTime]::FromOADate(44489)
$ws.Range("E260").Value = 16
$ws.Range("F260").Value = "Fruta"
$ws.Range("G260").Value = 100108
$ws.Range("H260").Value = "Tropicales y subtropicales"
$ws.Range("I260").Value = 100108006
$ws.Range("J260").Value = "Plátano"
$ws.Range("K260").Value = "Sin especificar"
$ws.Range("L260").Value = "Primera Pintón"
$ws.Range("M260").Value = 400
$ws.Range("N260").Value = 26000
$ws.Range("O260").Value = 27000
$ws.Range("P260").Value = 26500
$ws.Range("Q260").Value = "$/caja 20 kilos"
$ws.Range("R260").Value = "Ecuador"
$ws.Range("S260").Value = 1325
$ws.Range("T260").Value = 20

Write-Host "Done"
